{"js": "// Update the date in the \"Faite \u00e0 Bouchain, le ...\" line from 07/01/2021\n// to 21/02/2021, and drop the \" DK Digital\" company name that used to sit\n// between \"le profil de votre soci\u00e9t\u00e9\" and \" et son c\u0153ur d'activit\u00e9...\".\n\nconst body = context.document.body;\n\n// 1) Date: 07/01/2021 -> 21/02/2021\nconst dateResults = body.search(\"07/01/2021\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"21/02/2021\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Remove \" DK Digital\" (the leading space stays, coming from the run\n//    that used to hold \" et son c\u0153ur d'activit\u00e9...\").\nconst companyResults = body.search(\" DK Digital\", { matchCase: true });\ncompanyResults.load(\"items\");\nawait context.sync();\n\nif (companyResults.items.length > 0) {\n  companyResults.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date in the \"Faite \u00e0 Bouchain, le ...\" line from 07/01/2021\n# to 21/02/2021, and drop the \" DK Digital\" company name that used to sit\n# between \"le profil de votre soci\u00e9t\u00e9\" and \" et son c\u0153ur d'activit\u00e9...\".\n\n$d = $word.ActiveDocument\n\n# 1) Date: 07/01/2021 -> 21/02/2021\n$dateRange = $d.Content\n$dateFind = $dateRange.Find\n$dateFind.Text = \"07/01/2021\"\n$dateFind.MatchCase = $true\nif ($dateFind.Execute()) {\n    $dateRange.Text = \"21/02/2021\"\n}\n\n# 2) Remove \" DK Digital\" (the leading space stays, coming from the run\n#    that used to hold \" et son c\u0153ur d'activit\u00e9...\").\n$companyRange = $d.Content\n$companyFind = $companyRange.Find\n$companyFind.Text = \" DK Digital\"\n$companyFind.MatchCase = $true\nif ($companyFind.Execute()) {\n    $companyRange.Text = \"\"\n}\n"}
